# The new row 3 is an exact duplicate of row 2 (same match repeated in the
# source data). Copy row 2 -> row 3 instead of re-typing the values so the
# text-stored numbers (9, 13, 0, 0, 69.23) keep their "number stored as
# text" representation instead of being re-interpreted as numerics, and so
# any incidental characters (e.g. the non-breaking space after "Unadkat")
# are reproduced exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:K2").Copy()
$ws.Range("A3").PasteSpecial(-4163)
